# Add row 6 data for Insperity system to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 'Insperity'
$ws.Range("B6").Value = 'This system is primarily used for HR management and payroll processing. It supports employee onboarding, I-9 verifications, benefits elections, payroll data submission, time reporting, and vacation accrual tracking.'
$ws.Range("C6").Value = 'The client utilizes this system to manage onboarding processes, employee benefit elections, vacation time tracking, payroll processing, and manual compensation adjustments within employee profiles.'
$ws.Range("D6").Value = 'The HR team, led by Laurie Stewart (Senior Vice President of HR), oversees system administration and user access management.'
$ws.Range("E6").Value = 'HR inputs new hire details such as salary, start date, and manager into the system, triggering onboarding notifications and granting employee access to Insperity, while elevated access or role changes require an admin-submitted access request form reviewed and provisioned by the Insperity team.'
$ws.Range("F6").Value = 'Upon termination, HR updates the employee''s record with the termination date and reason. For standard employees, limited access is retained for viewing tax forms, while for elevated access users, an email is sent to Insperity to request access removal, and confirmation is documented via email.'
$ws.Range("G6").Value = 'Access is configured using a permission-based model, where specific permissions are selected through an access request form and provisioned by the Insperity team; users cannot create roles or grant access themselves.'
$ws.Range("H6").Value = 'No; Access is permission-based, and system administrators cannot create or modify roles. Changes to permissions require submitting an access request form, which is processed and provisioned by the Insperity team.'
$ws.Range("I6").Value = 'No; management intends to implement periodic reviews of elevated user access in the future, but no formal review process for all roles and permissions currently exists.'
$ws.Range("J6").Value = 'Yes; Users with privileged access include the HR team, Megan Hodgson, John, Laurie Stewart, and Heather Malcolm, who have admin access to manage system configurations and employee setups, with some actions requiring external intervention.'
$ws.Range("K6").Value = 'No; there are no interactive system accounts within the system.'
$ws.Range("L6").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("M6").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("N6").Value = 'Yes; but periodic user access reviews for elevated users are not currently being performed. Management plans to implement this process going forward.'
$ws.Range("O6").Value = 'No, the system does not have activity logging capabilities or audit trail functionality.'
$ws.Range("P6").Value = 'No; periodic activity reviews of user activity are not currently performed, but there are indirect controls in place such as payroll reviews.'
$ws.Range("Q6").Value = 'Users authenticate through direct login with enforced two-factor authentication for all accounts.'
$ws.Range("R6").Value = 'No; The client does not perform periodic reviews of the system''s authentication configurations.'
$ws.Range("S6").Value = 'Management does not have the ability to make any changes to this system. All system configurations, workflows, and code are managed entirely by the vendor.'
$ws.Range("T6").Value = 'Access to make changes is restricted to Megan Hodgson, John, and the HR team. Megan and John have administrative access due to the current size of the team and operational needs, while the HR team has admin access to manage application configurations.'
$ws.Range("U6").Value = 'Management does not have any separate environments for this system.'
$ws.Range("V6").Value = 'The client''s change management process includes the following steps: Access provisioning changes are requested via an access request form specifying permissions, approved by Insperity, and deployed by Insperity support without client-side provisioning. Payroll changes are requested via email, reviewed and approved by Megan and John, and implemented by Insperity after client approval. User access termination requests are sent via email, confirmed by Insperity, and access is disabled. Payroll data is reviewed semi-monthly by Megan and John, who approve and finalize it for submission. Changes are logged in a pay change history report, and future plans include implementing a new payroll system with formalized controls.'
$ws.Range("W6").Value = 'No, the vendor manages all updates, patches, and bug fixes directly in the production environment. The client has no ability to make changes or updates themselves, and there is no sandbox or QA environment for testing.'
$ws.Range("X6").Value = 'The system enforces segregation of duties through inherent controls. Users cannot independently provision access, make configuration changes, or modify payroll data, as these actions are handled by Insperity support. Additionally, the system includes audit logging and two-factor authentication to enhance security and monitoring.'
$ws.Range("Y6").Value = 'No periodic review of changes is performed.'
$ws.Range("Z6").Value = 'Yes, the system has automated jobs or interfaces with benefits administration sites and their website. These jobs perform functions related to benefits management and payroll journal entry preparation.'
$ws.Range("AA6").Value = 'Jobs are managed using vendor-managed tools provided by Insperity. These tools handle scheduling and monitoring of automated processes, while the organization does not currently use internal systems for these functions.'
$ws.Range("AB6").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AC6").Value = 'Data is stored in a vendor-managed database managed by Insperity.'
$ws.Range("AD6").Value = 'Backups are handled entirely by the vendor as part of their SaaS service.'
$ws.Range("AE6").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AF6").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AG6").Value = 'No, management does not perform regular SOC report reviews.'
